$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.7397403717041016
$ws.Range("E2").Value = 4432.100724521545
$ws.Range("G2").Value = 0.1006509244510039
$ws.Range("H2").Value = 0.1006509244510039
$ws.Range("I2").Value = 0.1006509244510039
$ws.Range("J2").Value = 0.0998012603099043
$ws.Range("K2").Value = 0.09644153875249714
$ws.Range("L2").Value = 0.09644153875249714
$ws.Range("M2").Value = 0.09644153875249714
$ws.Range("N2").Value = 0.09644153875249714
$ws.Range("O2").Value = 0.09644153875249714
$ws.Range("P2").Value = 0.09504701135541072
$ws.Range("Q2").Value = 0.09504701135541072
$ws.Range("R2").Value = 0.09504701135541072
$ws.Range("S2").Value = 0.09499665151736289
$ws.Range("T2").Value = 0.09499665151736289
$ws.Range("U2").Value = 0.09499665151736289
$ws.Range("V2").Value = 0.09470302653012871
$ws.Range("W2").Value = 0.09443318179132074
$ws.Range("X2").Value = 0.09439572562420165
$ws.Range("Y2").Value = 0.09439572562420165

# Row 3
$ws.Range("C3").Value = 0.8540639877319336
$ws.Range("E3").Value = 4419.639015056016
$ws.Range("G3").Value = 0.1010528132146395
$ws.Range("H3").Value = 0.1010528132146395
$ws.Range("I3").Value = 0.09834351336382258
$ws.Range("J3").Value = 0.097815640517698
$ws.Range("K3").Value = 0.097815640517698
$ws.Range("L3").Value = 0.097815640517698
$ws.Range("M3").Value = 0.097815640517698
$ws.Range("N3").Value = 0.097815640517698
$ws.Range("O3").Value = 0.09691305972859751
$ws.Range("P3").Value = 0.09691305972859751
$ws.Range("Q3").Value = 0.09691305972859751
$ws.Range("R3").Value = 0.09596795759226556
$ws.Range("S3").Value = 0.09470732940406171
$ws.Range("T3").Value = 0.09470732940406171
$ws.Range("U3").Value = 0.09470732940406171
$ws.Range("V3").Value = 0.09415760289339867
$ws.Range("W3").Value = 0.09415760289339867
$ws.Range("X3").Value = 0.09415760289339867
$ws.Range("Y3").Value = 0.09415280731103343

# Row 4
$ws.Range("C4").Value = 0.7959372997283936
$ws.Range("E4").Value = 4373.383029285585
$ws.Range("G4").Value = 0.09857589607701014
$ws.Range("H4").Value = 0.09438767143613172
$ws.Range("I4").Value = 0.09431317227244521
$ws.Range("J4").Value = 0.09431317227244521
$ws.Range("K4").Value = 0.09431317227244521
$ws.Range("L4").Value = 0.09431317227244521
$ws.Range("M4").Value = 0.09407006513978894
$ws.Range("N4").Value = 0.09407006513978894
$ws.Range("O4").Value = 0.09392110366297235
$ws.Range("P4").Value = 0.09339634563145457
$ws.Range("Q4").Value = 0.09339634563145457
$ws.Range("R4").Value = 0.09339634563145457
$ws.Range("S4").Value = 0.09339634563145457
$ws.Range("T4").Value = 0.09325113117515757
$ws.Range("U4").Value = 0.09325113117515757
$ws.Range("V4").Value = 0.09325113117515757
$ws.Range("W4").Value = 0.09325113117515757
$ws.Range("X4").Value = 0.09325113117515757
$ws.Range("Y4").Value = 0.09325113117515757

# Row 5
$ws.Range("C5").Value = 0.7499995231628418
$ws.Range("E5").Value = 4359.664404189797
$ws.Range("G5").Value = 0.1025310626375835
$ws.Range("H5").Value = 0.1006829219552243
$ws.Range("I5").Value = 0.09759714805102515
$ws.Range("J5").Value = 0.09698229653628856
$ws.Range("K5").Value = 0.09698229653628856
$ws.Range("L5").Value = 0.09333101260587309
$ws.Range("M5").Value = 0.09333101260587309
$ws.Range("N5").Value = 0.09333101260587309
$ws.Range("O5").Value = 0.09333101260587309
$ws.Range("P5").Value = 0.09333101260587309
$ws.Range("Q5").Value = 0.09333101260587309
$ws.Range("R5").Value = 0.09333101260587309
$ws.Range("S5").Value = 0.09326650462032929
$ws.Range("T5").Value = 0.09326650462032929
$ws.Range("U5").Value = 0.09326650462032929
$ws.Range("V5").Value = 0.09308273543073639
$ws.Range("W5").Value = 0.09308273543073639
$ws.Range("X5").Value = 0.09299423989880842
$ws.Range("Y5").Value = 0.09298371158264709

# Row 6
$ws.Range("C6").Value = 0.7499995231628418
$ws.Range("E6").Value = 4359.505261918309
$ws.Range("G6").Value = 0.1002683146796879
$ws.Range("H6").Value = 0.09746240666987097
$ws.Range("I6").Value = 0.09746240666987097
$ws.Range("J6").Value = 0.09746240666987097
$ws.Range("K6").Value = 0.09746240666987097
$ws.Range("L6").Value = 0.09620157530633339
$ws.Range("M6").Value = 0.09620157530633339
$ws.Range("N6").Value = 0.09620157530633339
$ws.Range("O6").Value = 0.09527995480605791
$ws.Range("P6").Value = 0.09527995480605791
$ws.Range("Q6").Value = 0.09304433697650311
$ws.Range("R6").Value = 0.09304433697650311
$ws.Range("S6").Value = 0.09304433697650311
$ws.Range("T6").Value = 0.09304433697650311
$ws.Range("U6").Value = 0.09304433697650311
$ws.Range("V6").Value = 0.09304433697650311
$ws.Range("W6").Value = 0.09304136648066039
$ws.Range("X6").Value = 0.09299342909772396
$ws.Range("Y6").Value = 0.09298060939411906

# Row 7
$ws.Range("C7").Value = 0.75
$ws.Range("E7").Value = 4355.198026534794
$ws.Range("G7").Value = 0.1027405855027473
$ws.Range("H7").Value = 0.09299508323591832
$ws.Range("I7").Value = 0.09299508323591832
$ws.Range("J7").Value = 0.09299508323591832
$ws.Range("K7").Value = 0.09299508323591832
$ws.Range("L7").Value = 0.09299508323591832
$ws.Range("M7").Value = 0.09299508323591832
$ws.Range("N7").Value = 0.09299508323591832
$ws.Range("O7").Value = 0.09299508323591832
$ws.Range("P7").Value = 0.09299508323591832
$ws.Range("Q7").Value = 0.09299508323591832
$ws.Range("R7").Value = 0.09289664769073672
$ws.Range("S7").Value = 0.09289664769073672
$ws.Range("T7").Value = 0.09289664769073672
$ws.Range("U7").Value = 0.09289664769073672
$ws.Range("V7").Value = 0.09289664769073672
$ws.Range("W7").Value = 0.09289664769073672
$ws.Range("X7").Value = 0.09289664769073672
$ws.Range("Y7").Value = 0.09289664769073672

# Row 8
$ws.Range("C8").Value = 0.7500007152557373
$ws.Range("E8").Value = 4366.22472262697
$ws.Range("G8").Value = 0.1008145362536919
$ws.Range("H8").Value = 0.09414456799122095
$ws.Range("I8").Value = 0.09414456799122095
$ws.Range("J8").Value = 0.09414456799122095
$ws.Range("K8").Value = 0.0937759268227244
$ws.Range("L8").Value = 0.0937759268227244
$ws.Range("M8").Value = 0.0937759268227244
$ws.Range("N8").Value = 0.0937759268227244
$ws.Range("O8").Value = 0.0934328585154385
$ws.Range("P8").Value = 0.09340251952258052
$ws.Range("Q8").Value = 0.09340251952258052
$ws.Range("R8").Value = 0.09340251952258052
$ws.Range("S8").Value = 0.09340251952258052
$ws.Range("T8").Value = 0.09340251952258052
$ws.Range("U8").Value = 0.09328810381095509
$ws.Range("V8").Value = 0.09321936165036326
$ws.Range("W8").Value = 0.09321936165036326
$ws.Range("X8").Value = 0.09321936165036326
$ws.Range("Y8").Value = 0.09311159303366412

# Row 9
$ws.Range("C9").Value = 0.7187490463256836
$ws.Range("E9").Value = 4415.032113387887
$ws.Range("G9").Value = 0.09935604244664364
$ws.Range("H9").Value = 0.09935604244664364
$ws.Range("I9").Value = 0.09910093414134508
$ws.Range("J9").Value = 0.09782516742229239
$ws.Range("K9").Value = 0.09765355444904092
$ws.Range("L9").Value = 0.09486511341133677
$ws.Range("M9").Value = 0.09444921745876
$ws.Range("N9").Value = 0.09444921745876
$ws.Range("O9").Value = 0.09444921745876
$ws.Range("P9").Value = 0.09444921745876
$ws.Range("Q9").Value = 0.09444921745876
$ws.Range("R9").Value = 0.09444921745876
$ws.Range("S9").Value = 0.09444921745876
$ws.Range("T9").Value = 0.09444921745876
$ws.Range("U9").Value = 0.09444921745876
$ws.Range("V9").Value = 0.09438434236543031
$ws.Range("W9").Value = 0.0943112291509059
$ws.Range("X9").Value = 0.09406365976125713
$ws.Range("Y9").Value = 0.09406300415960792

# Row 10
$ws.Range("C10").Value = 0.7656230926513672
$ws.Range("E10").Value = 4374.698399875145
$ws.Range("G10").Value = 0.1026872366763494
$ws.Range("H10").Value = 0.09448790280966815
$ws.Range("I10").Value = 0.09448790280966815
$ws.Range("J10").Value = 0.09338716074647138
$ws.Range("K10").Value = 0.09338716074647138
$ws.Range("L10").Value = 0.09338716074647138
$ws.Range("M10").Value = 0.09338716074647138
$ws.Range("N10").Value = 0.09338716074647138
$ws.Range("O10").Value = 0.09338716074647138
$ws.Range("P10").Value = 0.09338716074647138
$ws.Range("Q10").Value = 0.09338716074647138
$ws.Range("R10").Value = 0.09333675714086492
$ws.Range("S10").Value = 0.09333675714086492
$ws.Range("T10").Value = 0.09333675714086492
$ws.Range("U10").Value = 0.09333675714086492
$ws.Range("V10").Value = 0.09327677192739073
$ws.Range("W10").Value = 0.09327677192739073
$ws.Range("X10").Value = 0.09327677192739073
$ws.Range("Y10").Value = 0.09327677192739073

# Row 11
$ws.Range("C11").Value = 0.7500112056732178
$ws.Range("E11").Value = 4384.671138216754
$ws.Range("G11").Value = 0.1023664668111599
$ws.Range("H11").Value = 0.09362718720498697
$ws.Range("I11").Value = 0.09362718720498697
$ws.Range("J11").Value = 0.09362718720498697
$ws.Range("K11").Value = 0.09362718720498697
$ws.Range("L11").Value = 0.09362718720498697
$ws.Range("M11").Value = 0.09362718720498697
$ws.Range("N11").Value = 0.09362718720498697
$ws.Range("O11").Value = 0.09362718720498697
$ws.Range("P11").Value = 0.0935271225131821
$ws.Range("Q11").Value = 0.0935271225131821
$ws.Range("R11").Value = 0.0935271225131821
$ws.Range("S11").Value = 0.0935271225131821
$ws.Range("T11").Value = 0.0935271225131821
$ws.Range("U11").Value = 0.0935271225131821
$ws.Range("V11").Value = 0.09347117228492696
$ws.Range("W11").Value = 0.09347117228492696
$ws.Range("X11").Value = 0.09347117228492696
$ws.Range("Y11").Value = 0.09347117228492696

